# example with comparison of calculations
# example of detailed calculations of s2 without using VAR.S function
# (to be deleted after checked by user) -- #23

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example")
$ws.Activate()

# --- labels, in the same order the original author typed them so the
#     shared-string table comes out in the same sequence -------------------
$ws.Range("R35").Value = "squared errors"
$ws.Range("R54").Value = "comparison R33 vs R52"
$ws.Range("R51").Value = "sum of squared errors / 13"

# fill style (green) must exist before any quote-prefixed text style below,
# to line up with the order the workbook's cellXfs table was authored in
$ws.Range("S54").Interior.Color = 5296274

# --- squared-errors column (R36:R49 = R3^2 .. R16^2) -----------------------
$ws.Range("R36").Formula = "=R3^2"
$ws.Range("S36").Value = "'-> R3^2"

$ws.Range("R37").Formula = "=R4^2"
$ws.Range("R38").Formula = "=R5^2"
$ws.Range("R39").Formula = "=R6^2"
$ws.Range("R40").Formula = "=R7^2"
$ws.Range("R41").Formula = "=R8^2"
$ws.Range("R42").Formula = "=R9^2"
$ws.Range("R43").Formula = "=R10^2"
$ws.Range("R44").Formula = "=R11^2"
$ws.Range("R45").Formula = "=R12^2"
$ws.Range("R46").Formula = "=R13^2"
$ws.Range("R47").Formula = "=R14^2"
$ws.Range("R48").Formula = "=R15^2"

# --- sum of squared errors / (n-1) -----------------------------------------
$ws.Range("S51").Formula = "=SUM(R36:R49)/13"
$ws.Range("T51").Value = "'-> =SUM(R36:R49)/13"

$ws.Range("R49").Formula = "=R16^2"
$ws.Range("S49").Value = "'-> R16^2"

# --- comparison against VAR.S result in R33 ---------------------------------
$ws.Range("S54").Formula = "=S51=R33"
$ws.Range("T54").Value = "'-> =S51=R33"

# --- hide the per-species detail columns, widen the new R column -----------
$ws.Range("G1:M1").EntireColumn.Hidden = $true
$ws.Columns("R").ColumnWidth = 21.6

# --- arrow connector pointing from the squared-errors block to the sum -----
$cxn = $ws.Shapes.AddConnector(1, 8092545, 1641523, 8092545, 3728364)
$cxn.Name = "Straight Arrow Connector 2"
$cxn.Line.EndArrowheadStyle = 2

# --- view state: zoom out, scroll back to the top, select the new cell -----
$excel.ActiveWindow.Zoom = 55
$ws.Range("T54").Select()
